$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("E2").Value = 77
$ws.Range("G2").Value = 94
$ws.Range("I2").Value = 121
$ws.Range("D3").Value = 142
$ws.Range("F3").Value = 148
$ws.Range("G3").Value = 150
$ws.Range("H3").Value = 163
$ws.Range("I3").Value = 202
$ws.Range("C9").Value = 510
$ws.Range("E9").Value = 519
$ws.Range("F9").Value = 591
$ws.Range("H9").Value = 484
$ws.Range("I9").Value = 520
$ws.Range("B10").Value = 1429
$ws.Range("C10").Value = 1676
$ws.Range("D10").Value = 1896
$ws.Range("E10").Value = 2341
$ws.Range("F10").Value = 2230
$ws.Range("G10").Value = 930
$ws.Range("H10").Value = 642
$ws.Range("B11").Value = 1968
$ws.Range("C11").Value = 2352
$ws.Range("D11").Value = 2592
$ws.Range("E11").Value = 3103
$ws.Range("F11").Value = 3080
$ws.Range("G11").Value = 1634
$ws.Range("H11").Value = 1422
$ws.Range("I11").Value = 1756

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I14").Value = 7
$ws.Range("C18").Value = 4
$ws.Range("E21").Value = 25
$ws.Range("F23").Value = 19
$ws.Range("F27").Value = 29
$ws.Range("D28").Value = 107
$ws.Range("C32").Value = 103
$ws.Range("D32").Value = 101
$ws.Range("E32").Value = 178
$ws.Range("G36").Value = 47
$ws.Range("I36").Value = 77
$ws.Range("E46").Value = 9
$ws.Range("D47").Value = 57
$ws.Range("F50").Value = 62
$ws.Range("B53").Value = 282
$ws.Range("C53").Value = 416
$ws.Range("E53").Value = 797
$ws.Range("F53").Value = 660
$ws.Range("H53").Value = 234
$ws.Range("I53").Value = 323
$ws.Range("H59").Value = 4
$ws.Range("E65").Value = 68
$ws.Range("B70").Value = 27
$ws.Range("G70").Value = 39
$ws.Range("E71").Value = 9
$ws.Range("G72").Value = 8
$ws.Range("F74").Value = 93
$ws.Range("I78").Value = 21
$ws.Range("D87").Value = 19
$ws.Range("E92").Value = 43
$ws.Range("H92").Value = 25
$ws.Range("H95").Value = 20
$ws.Range("B96").Value = 18
$ws.Range("G96").Value = 8
$ws.Range("B99").Value = 1968
$ws.Range("C99").Value = 2352
$ws.Range("D99").Value = 2592
$ws.Range("E99").Value = 3103
$ws.Range("F99").Value = 3080
$ws.Range("G99").Value = 1634
$ws.Range("H99").Value = 1422
$ws.Range("I99").Value = 1756

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("F7").Value = 29

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("E7").Value = 4
$ws.Range("E9").Value = 25

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("C8").Value = 63
$ws.Range("D8").Value = 50
$ws.Range("E8").Value = 109
$ws.Range("C9").Value = 103
$ws.Range("D9").Value = 101
$ws.Range("E9").Value = 178

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 4
$ws.Range("G3").Value = 9
$ws.Range("G9").Value = 47
$ws.Range("I9").Value = 77

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("E2").Value = 5
$ws.Range("I3").Value = 32
$ws.Range("F7").Value = 69
$ws.Range("H7").Value = 80
$ws.Range("B8").Value = 233
$ws.Range("C8").Value = 354
$ws.Range("E8").Value = 705
$ws.Range("F8").Value = 568
$ws.Range("H8").Value = 116
$ws.Range("B9").Value = 282
$ws.Range("C9").Value = 416
$ws.Range("E9").Value = 797
$ws.Range("F9").Value = 660
$ws.Range("H9").Value = 234
$ws.Range("I9").Value = 323

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("B6").Value = 24
$ws.Range("G6").Value = 24
$ws.Range("B7").Value = 27
$ws.Range("G7").Value = 39

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("F6").Value = 34
$ws.Range("F7").Value = 62

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("E6").Value = 21
$ws.Range("E8").Value = 68

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("D3").Value = 2

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("D9").Value = 19

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I4").Value = 9
$ws.Range("I6").Value = 21

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("D3").Value = 18
$ws.Range("D9").Value = 107

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("D7").Value = 43
$ws.Range("D8").Value = 57

$ws = $wb.Worksheets.Item('River North')
$ws.Range("F5").Value = 12
$ws.Range("F7").Value = 93

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("H3").Value = 3

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("E7").Value = 8
$ws.Range("E9").Value = 43
$ws.Range("H9").Value = 25

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("F6").Value = 19

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("C4").Value = 3

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("C6").Value = 4

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("F3").Value = 3
$ws.Range("F6").Value = 19

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I3").Value = 1

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I6").Value = 7

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 4

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("H5").Value = 12
$ws.Range("H6").Value = 20

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("G2").Value = 1
$ws.Range("B6").Value = 18
$ws.Range("B7").Value = 18
$ws.Range("G7").Value = 8

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 8

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("E6").Value = 8

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("E7").Value = 9

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("E2").Value = 2

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("E7").Value = 9
